# Deploying to gh-pages from @ LinuxForHealth/alvearie-fhir-ig@80fa500adfae01c9a5dd7ef65e90accc96781b5c
# - rebrand IBM -> LinuxForHealth (URL + Publisher), bump Version/Date,
#   and move the ele-1/ext-1 constraint text off the base "Extension" row
#   (it belongs on "Extension.extension" where it already also appears).

$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/employee-sub-business-unit"
$wsMeta.Range("B3").Value = "8.0.0"
$wsMeta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$wsMeta.Range("B9").Value = "LinuxForHealth Team"

$wsElements = $wb.Worksheets.Item("Elements")
# Fixed Value on the Extension.url row also carries the old IBM URL.
$wsElements.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/employee-sub-business-unit"
# The combined ele-1/ext-1 invariant no longer belongs on the base "Extension" row.
$wsElements.Range("AI2").Value = ""
